$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh timestamp + account numbers (quote-prefixed, forces text) + folder date bump ---
$ws.Range("A2").Value = 45230.625
$ws.Range("B2").Value = "'127520840434805,103482326003878,17841461742288388,17841456036806884"
$ws.Range("C2").Value = "/Users/jishuliu/Desktop/hkshop/data/20231031_superdelivery/【日本直送】 ササガワ　ＩＴ事業部（すべてのジャンル）  荷札シール　取扱注意  行李標籤貼紙 小心輕放  一色入"
$ws.Range("D2").Value = "https://shopage.s3.amazonaws.com/media/f854/615273998674_84835830603419215654.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_46914841642146126210.jpg"

# --- Row 3: refresh timestamp + same account numbers + new folder (HEIKO item) + 3 image urls + trailing empty cell ---
$ws.Range("A3").Value = 45230.6458333333
$ws.Range("B3").Value = "'127520840434805,103482326003878,17841461742288388,17841456036806884"
$ws.Range("C3").Value = "/Users/jishuliu/Desktop/hkshop/data/20231031_superdelivery/【日本直送】 HEIKO(ヘイコー)  ヘイコー 注意喚起シール カッター注意 48枚  Heiko 警告貼紙切割器警告 48 張  一色入"
$ws.Range("D3").Value = "https://shopage.s3.amazonaws.com/media/f854/615273998674_32746879691743666134.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_78459562396381333295.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_04964299814229126722.jpg"

# E3 needs to exist as an (empty-string) shared-string cell, not an empty/cleared cell -
# enter it quote-prefixed then strip the resulting quote-prefix formatting back off.
$ws.Range("E3").Value = "'"
$ws.Range("E3").ClearFormats()

# --- Row 4 no longer exists: drop it ---
$ws.Rows(4).Delete()

# --- column widths (closest reachable to the author's manual resize) ---
$ws.Columns(1).ColumnWidth = 35.42857142857143
$ws.Columns(3).ColumnWidth = 192.71428571428572
$ws.Columns(4).ColumnWidth = 366.7142857142857

# --- cursor moved on to where the next row would go ---
$ws.Range("B11").Select()
